$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Duplicate the "Optimal_length" column (C) into a new column D before
# renaming anything, so both columns end up referencing the same shared
# string values.
$ws.Range("C1:C14").Copy($ws.Range("D1:D14"))

# Rename column headers.
$ws.Range("C1").Value = "Optimal_length_upravene"
$ws.Range("D1").Value = "Optimal_length_WU"

# Apply number format 0.00 to F0M and both optimal length columns.
$ws.Range("B2:D14").NumberFormat = "0.00"

# Adjust column widths (values chosen so the stored OOXML width rounds to
# the authored 22.5703125 / 32.42578125 given this engine's width rounding).
$ws.Columns.Item(3).ColumnWidth = 21.7
$ws.Columns.Item(4).ColumnWidth = 31.7

# Update selection to match authored state.
$ws.Range("C7").Select()
